# Generate Report for Handoff
# Updates the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps for the
# d68ca6e4-a0cc-4175-8309-ddb4a376cf28.md file row, reflecting a newer handoff.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 5 is the d68ca6e4-... file, column D = "Latest Handoff Date"
$wsOverview.Range("D5").Value = "2016-03-24 02:11:04"

# zh-cn detail sheet: row 5 is the d68ca6e4-... file, column E = "Latest Handoff Datetime"
$wsZhCn.Range("E5").Value = "2016-03-24 02:10:55"

# de-de detail sheet: row 5 is the d68ca6e4-... file, column E = "Latest Handoff Datetime"
$wsDeDe.Range("E5").Value = "2016-03-24 02:11:04"
